$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 14.98669760470756
$ws.Cells.Item(2, 3).Value = 12.1461167406979
$ws.Cells.Item(2, 4).Value = 4.80754083459285
$ws.Cells.Item(2, 6).Value = 23.74236467974593
$ws.Cells.Item(2, 7).Value = 3.62712486080466
$ws.Cells.Item(2, 9).Value = 21.92968259560744
$ws.Cells.Item(2, 12).Value = 10.93843459763619
$ws.Cells.Item(2, 13).Value = 14.93099595568168
$ws.Cells.Item(2, 15).Value = 21.25138843000688

$ws.Cells.Item(3, 2).Value = 14.37999412778179
$ws.Cells.Item(3, 3).Value = 11.89522532001232
$ws.Cells.Item(3, 4).Value = 4.76033704974826
$ws.Cells.Item(3, 6).Value = 23.78763034473561
$ws.Cells.Item(3, 7).Value = 3.629193155550179
$ws.Cells.Item(3, 9).Value = 22.07615201440281
$ws.Cells.Item(3, 12).Value = 10.96069001648453
$ws.Cells.Item(3, 13).Value = 14.80827243883534
$ws.Cells.Item(3, 15).Value = 21.34241164285425

$ws.Cells.Item(4, 2).Value = 13.994729247701
$ws.Cells.Item(4, 3).Value = 11.73778610231282
$ws.Cells.Item(4, 4).Value = 4.730891512314495
$ws.Cells.Item(4, 6).Value = 23.82371044960678
$ws.Cells.Item(4, 7).Value = 3.630530626851026
$ws.Cells.Item(4, 9).Value = 22.17156813704139
$ws.Cells.Item(4, 12).Value = 10.97611524154257
$ws.Cells.Item(4, 13).Value = 14.73422723323268
$ws.Cells.Item(4, 15).Value = 21.4045383191226

$ws.Cells.Item(5, 2).Value = 13.83473663620841
$ws.Cells.Item(5, 3).Value = 11.67283842474492
$ws.Cells.Item(5, 4).Value = 4.718782546741959
$ws.Cells.Item(5, 6).Value = 23.84048893074279
$ws.Cells.Item(5, 7).Value = 3.631092691971704
$ws.Cells.Item(5, 9).Value = 22.21182940755471
$ws.Cells.Item(5, 12).Value = 10.98284372522348
$ws.Cells.Item(5, 13).Value = 14.70440682564251
$ws.Cells.Item(5, 15).Value = 21.43141765470623

$ws.Cells.Item(6, 2).Value = 13.80799583692217
$ws.Cells.Item(6, 3).Value = 11.66200809412918
$ws.Cells.Item(6, 4).Value = 4.716765431378326
$ws.Cells.Item(6, 6).Value = 23.84340007931339
$ws.Cells.Item(6, 7).Value = 3.631187052892714
$ws.Cells.Item(6, 9).Value = 22.21859800464722
$ws.Cells.Item(6, 12).Value = 10.98398771277338
$ws.Cells.Item(6, 13).Value = 14.69947724081364
$ws.Cells.Item(6, 15).Value = 21.4359751112418

$ws.Cells.Item(7, 2).Value = 13.99258335597062
$ws.Cells.Item(7, 3).Value = 11.73691330851759
$ws.Cells.Item(7, 4).Value = 4.730728641916641
$ws.Cells.Item(7, 6).Value = 23.82392833803017
$ws.Cells.Item(7, 7).Value = 3.630538138013778
$ws.Cells.Item(7, 9).Value = 22.17210553389752
$ws.Cells.Item(7, 12).Value = 10.97620419235633
$ws.Cells.Item(7, 13).Value = 14.7338236006709
$ws.Cells.Item(7, 15).Value = 21.4048945065164

$ws.Cells.Item(8, 2).Value = 14.78027052437946
$ws.Cells.Item(8, 3).Value = 12.06034986948672
$ws.Cells.Item(8, 4).Value = 4.791364814316936
$ws.Cells.Item(8, 6).Value = 23.7562477481201
$ws.Cells.Item(8, 7).Value = 3.627824025398963
$ws.Cells.Item(8, 9).Value = 21.97904693663778
$ws.Cells.Item(8, 12).Value = 10.94574288750539
$ws.Cells.Item(8, 13).Value = 14.88842457341349
$ws.Cells.Item(8, 15).Value = 21.28147472954315

$ws.Cells.Item(9, 2).Value = 16.2157781291287
$ws.Cells.Item(9, 3).Value = 12.66510357647777
$ws.Cells.Item(9, 4).Value = 4.906321812693662
$ws.Cells.Item(9, 6).Value = 23.68959316426503
$ws.Cells.Item(9, 7).Value = 3.623035024891098
$ws.Cells.Item(9, 9).Value = 21.64398568666124
$ws.Cells.Item(9, 12).Value = 10.89997747415908
$ws.Cells.Item(9, 13).Value = 15.20085728030844
$ws.Cells.Item(9, 15).Value = 21.08922370377971

$ws.Cells.Item(10, 2).Value = 17.19480520424838
$ws.Cells.Item(10, 3).Value = 13.08818058555424
$ws.Cells.Item(10, 4).Value = 4.98800234477722
$ws.Cells.Item(10, 6).Value = 23.68124770725855
$ws.Cells.Item(10, 7).Value = 3.619838263633759
$ws.Cells.Item(10, 9).Value = 21.42437775706716
$ws.Cells.Item(10, 12).Value = 10.87487051809638
$ws.Cells.Item(10, 13).Value = 15.43452970787552
$ws.Cells.Item(10, 15).Value = 20.97869026157171

$ws.Cells.Item(11, 2).Value = 17.62225414212185
$ws.Cells.Item(11, 3).Value = 13.27541850395603
$ws.Cells.Item(11, 4).Value = 5.02448150294159
$ws.Cells.Item(11, 6).Value = 23.68631964724275
$ws.Cells.Item(11, 7).Value = 3.618453100739098
$ws.Cells.Item(11, 9).Value = 21.33024813405085
$ws.Cells.Item(11, 12).Value = 10.86529733978219
$ws.Cells.Item(11, 13).Value = 15.54140212865589
$ws.Cells.Item(11, 15).Value = 20.93514960349458

$ws.Cells.Item(12, 2).Value = 17.78143926230084
$ws.Cells.Item(12, 3).Value = 13.34552257189124
$ws.Cells.Item(12, 4).Value = 5.038191510363253
$ws.Cells.Item(12, 6).Value = 23.68951705984742
$ws.Cells.Item(12, 7).Value = 3.61793845088873
$ws.Cells.Item(12, 9).Value = 21.29543468225157
$ws.Cells.Item(12, 12).Value = 10.86193782017288
$ws.Cells.Item(12, 13).Value = 15.58192676141179
$ws.Cells.Item(12, 15).Value = 20.91963710395504

$ws.Cells.Item(13, 2).Value = 17.7472766128052
$ws.Cells.Item(13, 3).Value = 13.33046066298972
$ws.Cells.Item(13, 4).Value = 5.035243534212722
$ws.Cells.Item(13, 6).Value = 23.6887716460835
$ws.Cells.Item(13, 7).Value = 3.61804885125301
$ws.Cells.Item(13, 9).Value = 21.30289538002253
$ws.Cells.Item(13, 12).Value = 10.86264954115651
$ws.Cells.Item(13, 13).Value = 15.57319708515505
$ws.Cells.Item(13, 15).Value = 20.92293452480112

$ws.Cells.Item(14, 2).Value = 17.63540466198765
$ws.Cells.Item(14, 3).Value = 13.28120228063609
$ws.Cells.Item(14, 4).Value = 5.025611551044715
$ws.Cells.Item(14, 6).Value = 23.68655710823221
$ws.Cells.Item(14, 7).Value = 3.61841056243277
$ws.Cells.Item(14, 9).Value = 21.32736733891524
$ws.Cells.Item(14, 12).Value = 10.86501562724767
$ws.Cells.Item(14, 13).Value = 15.54473518564975
$ws.Cells.Item(14, 15).Value = 20.93385379906071

$ws.Cells.Item(15, 2).Value = 17.56652794809292
$ws.Cells.Item(15, 3).Value = 13.25092469909409
$ws.Cells.Item(15, 4).Value = 5.019697970714833
$ws.Cells.Item(15, 6).Value = 23.68536693439004
$ws.Cells.Item(15, 7).Value = 3.61863340634709
$ws.Cells.Item(15, 9).Value = 21.34246543950268
$ws.Cells.Item(15, 12).Value = 10.86649951230816
$ws.Cells.Item(15, 13).Value = 15.52730769783174
$ws.Cells.Item(15, 15).Value = 20.94066936866649

$ws.Cells.Item(16, 2).Value = 17.16650063277617
$ws.Cells.Item(16, 3).Value = 13.07583488900332
$ws.Cells.Item(16, 4).Value = 4.985604157382084
$ws.Cells.Item(16, 6).Value = 23.68109486118056
$ws.Cells.Item(16, 7).Value = 3.619930172711224
$ws.Cells.Item(16, 9).Value = 21.43064560538504
$ws.Cells.Item(16, 12).Value = 10.87553332745335
$ws.Cells.Item(16, 13).Value = 15.42755448077551
$ws.Cells.Item(16, 15).Value = 20.98167192035707

$ws.Cells.Item(17, 2).Value = 16.91642901262459
$ws.Cells.Item(17, 3).Value = 12.96705026470391
$ws.Cells.Item(17, 4).Value = 4.964510518283812
$ws.Cells.Item(17, 6).Value = 23.68074702228726
$ws.Cells.Item(17, 7).Value = 3.620743349451449
$ws.Cells.Item(17, 9).Value = 21.48622043912118
$ws.Cells.Item(17, 12).Value = 10.88154855696911
$ws.Cells.Item(17, 13).Value = 15.36648600201232
$ws.Cells.Item(17, 15).Value = 21.00855666993507

$ws.Cells.Item(18, 2).Value = 16.77091331969661
$ws.Cells.Item(18, 3).Value = 12.90399135736551
$ws.Cells.Item(18, 4).Value = 4.952314597074976
$ws.Cells.Item(18, 6).Value = 23.68138168456149
$ws.Cells.Item(18, 7).Value = 3.621217570375151
$ws.Cells.Item(18, 9).Value = 21.5187287204813
$ws.Cells.Item(18, 12).Value = 10.88518232077237
$ws.Cells.Item(18, 13).Value = 15.33141670414955
$ws.Cells.Item(18, 15).Value = 21.02465427356887

$ws.Cells.Item(19, 2).Value = 16.7213589108288
$ws.Cells.Item(19, 3).Value = 12.88255824979096
$ws.Cells.Item(19, 4).Value = 4.948174573399917
$ws.Cells.Item(19, 6).Value = 23.68173986566377
$ws.Cells.Item(19, 7).Value = 3.621379251869846
$ws.Cells.Item(19, 9).Value = 21.52982870122559
$ws.Cells.Item(19, 12).Value = 10.88644253174879
$ws.Cells.Item(19, 13).Value = 15.31955326341106
$ws.Cells.Item(19, 15).Value = 21.03021338471662

$ws.Cells.Item(20, 2).Value = 16.94322435721062
$ws.Cells.Item(20, 3).Value = 12.97868152823581
$ws.Cells.Item(20, 4).Value = 4.96676258258569
$ws.Cells.Item(20, 6).Value = 23.6806976432031
$ws.Cells.Item(20, 7).Value = 3.620656112704852
$ws.Cells.Item(20, 9).Value = 21.48024818310899
$ws.Cells.Item(20, 12).Value = 10.88089022169512
$ws.Cells.Item(20, 13).Value = 15.37298128764898
$ws.Cells.Item(20, 15).Value = 21.00562906025564

$ws.Cells.Item(21, 2).Value = 17.66833766906998
$ws.Cells.Item(21, 3).Value = 13.29569270755113
$ws.Cells.Item(21, 4).Value = 5.028443566845856
$ws.Cells.Item(21, 6).Value = 23.68717291538702
$ws.Cells.Item(21, 7).Value = 3.618304051307999
$ws.Cells.Item(21, 9).Value = 21.32015675487769
$ws.Cells.Item(21, 12).Value = 10.86431344238512
$ws.Cells.Item(21, 13).Value = 15.55309388349256
$ws.Cells.Item(21, 15).Value = 20.93062002512978

$ws.Cells.Item(22, 2).Value = 18.12657916237032
$ws.Cells.Item(22, 3).Value = 13.49820162608233
$ws.Cells.Item(22, 4).Value = 5.068147443337899
$ws.Cells.Item(22, 6).Value = 23.69884659270526
$ws.Cells.Item(22, 7).Value = 3.616824416175692
$ws.Cells.Item(22, 9).Value = 21.22037454732215
$ws.Cells.Item(22, 12).Value = 10.85502778748125
$ws.Cells.Item(22, 13).Value = 15.67111183261281
$ws.Cells.Item(22, 15).Value = 20.88728525054037

$ws.Cells.Item(23, 2).Value = 17.88347068216857
$ws.Cells.Item(23, 3).Value = 13.39056162474543
$ws.Cells.Item(23, 4).Value = 5.047014459546729
$ws.Cells.Item(23, 6).Value = 23.6919350906035
$ws.Cells.Item(23, 7).Value = 3.617608873601692
$ws.Cells.Item(23, 9).Value = 21.2731861727388
$ws.Cells.Item(23, 12).Value = 10.85984210785056
$ws.Cells.Item(23, 13).Value = 15.60810487390324
$ws.Cells.Item(23, 15).Value = 20.90989148456788

$ws.Cells.Item(24, 2).Value = 16.93111561284401
$ws.Cells.Item(24, 3).Value = 12.97342464124381
$ws.Cells.Item(24, 4).Value = 4.965744638450451
$ws.Cells.Item(24, 6).Value = 23.68071736775888
$ws.Cells.Item(24, 7).Value = 3.620695531517833
$ws.Cells.Item(24, 9).Value = 21.48294650333319
$ws.Cells.Item(24, 12).Value = 10.88118730831085
$ws.Cells.Item(24, 13).Value = 15.37004464293681
$ws.Cells.Item(24, 15).Value = 21.00695063571275

$ws.Cells.Item(25, 2).Value = 15.84014303224503
$ws.Cells.Item(25, 3).Value = 12.50500887156499
$ws.Cells.Item(25, 4).Value = 4.875684457958293
$ws.Cells.Item(25, 6).Value = 23.70051079681707
$ws.Cells.Item(25, 7).Value = 3.624273832581747
$ws.Cells.Item(25, 9).Value = 21.72996654075148
$ws.Cells.Item(25, 12).Value = 10.91086244154512
$ws.Cells.Item(25, 13).Value = 15.20085728030844
$ws.Cells.Item(25, 15).Value = 21.13586499829437
